# Update agency_data sheet: rename several "Risk Management" program labels
# to their new, more formal "Framework" names.
#
# Mapping (old -> new):
#   Water Resource Management   -> National Water Resource Management Sector Framework
#   Flood Risk Management       -> Overarching Flood Risk Management Framework
#   Drought Risk Management     -> Overarching National Drought Risk Management Framework
#   Disaster Risk Management    -> National Disaster Risk Management Sector Framework
#
# These labels live in column B ("Program") of the agency_data worksheet, on
# rows 2, 3, 4, 18, 19 and 32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("agency_data")

$ws.Range("B2").Value  = "National Water Resource Management Sector Framework"
$ws.Range("B3").Value  = "Overarching Flood Risk Management Framework"
$ws.Range("B4").Value  = "Overarching National Drought Risk Management Framework"
$ws.Range("B18").Value = "National Disaster Risk Management Sector Framework"
$ws.Range("B19").Value = "Overarching Flood Risk Management Framework"
$ws.Range("B32").Value = "Overarching National Drought Risk Management Framework"

# Keep the selection/active cell consistent with the edited workbook.
$ws.Range("B18").Select() | Out-Null
